$wb = $excel.ActiveWorkbook

# --- survey sheet ---
$survey = $wb.Worksheets.Item("survey")

# Remove the "Date Entered" / today() / readonly row (old row 22).
$survey.Rows.Item(22).Delete()

# Insert a new blank row before "Date of Appointment" (now row 25 after the
# deletion above) to hold the new "Date of Expected Blood Draw" question.
$survey.Rows.Item(25).Insert()
$survey.Range("A25").Value = "date"
$survey.Range("B25").Value = "date_appoint1"
$survey.Range("C25").Value = "Date of Expected Blood Draw"

# Rename the "Lab test:" label to "Blood Draw:" (row 24 after the delete).
$survey.Range("C24").Value = "Blood Draw:"

# --- choices sheet ---
$choices = $wb.Worksheets.Item("choices")

# Rename the "lab test appointment" choice to "blood draw appointment".
$choices.Range("B5").Value = "blood draw appointment"
$choices.Range("C5").Value = "Blood Draw Appointment
"
# Setting a multi-line value auto-expands the row; restore its natural
# (default) height so the row doesn't pick up a stray explicit height.
$choices.Rows.Item(5).AutoFit()
